$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new display text.
# Price-column (D) values that look numeric are written with a leading
# apostrophe so Excel keeps them as text (matching the workbook's original
# inlineStr/text formatting) instead of silently converting them to numbers.
$updates = @(
    @{ Cell = 'D2'; Value = '30.140.27' }
    @{ Cell = 'E2'; Value = '  +0.63%  ' }
    @{ Cell = 'D3'; Value = '1.921.42' }
    @{ Cell = 'E3'; Value = '  +2.80%  ' }
    @{ Cell = 'E4'; Value = '  +0.16%  ' }
    @{ Cell = 'D5'; Value = '''319.56' }
    @{ Cell = 'E5'; Value = '  +0.14%  ' }
    @{ Cell = 'E6'; Value = '  +0.12%  ' }
    @{ Cell = 'D7'; Value = '''0.5072' }
    @{ Cell = 'E7'; Value = '  -0.37%  ' }
    @{ Cell = 'D8'; Value = '''0.4079' }
    @{ Cell = 'E8'; Value = '  +3.55%  ' }
    @{ Cell = 'D9'; Value = '''0.08338' }
    @{ Cell = 'E9'; Value = '  +1.65%  ' }
    @{ Cell = 'D10'; Value = '''1.118' }
    @{ Cell = 'E10'; Value = '  +2.24%  ' }
    @{ Cell = 'E11'; Value = '  -0.22%  ' }
    @{ Cell = 'D12'; Value = '''24.08' }
    @{ Cell = 'E12'; Value = '  +4.91%  ' }
    @{ Cell = 'D13'; Value = '1.922.13' }
    @{ Cell = 'E13'; Value = '  +2.52%  ' }
    @{ Cell = 'D14'; Value = '''6.429' }
    @{ Cell = 'E14'; Value = '  +2.33%  ' }
    @{ Cell = 'D15'; Value = '''7.254' }
    @{ Cell = 'E15'; Value = '  +0.91%  ' }
    @{ Cell = 'D16'; Value = '''1.002' }
    @{ Cell = 'E16'; Value = '  +0.06%  ' }
    @{ Cell = 'D17'; Value = '''92.70' }
    @{ Cell = 'E17'; Value = '  +0.74%  ' }
    @{ Cell = 'E18'; Value = '  +0.84%  ' }
    @{ Cell = 'D19'; Value = '''0.06516' }
    @{ Cell = 'E19'; Value = '  +2.00%  ' }
    @{ Cell = 'D20'; Value = '''18.51' }
    @{ Cell = 'E20'; Value = '  +3.49%  ' }
    @{ Cell = 'D21'; Value = '''1.002' }
    @{ Cell = 'E21'; Value = '  +0.11%  ' }
    @{ Cell = 'D22'; Value = '''5.959' }
    @{ Cell = 'E22'; Value = '  +2.28%  ' }
    @{ Cell = 'D23'; Value = '30.147.63' }
    @{ Cell = 'E23'; Value = '  +0.69%  ' }
    @{ Cell = 'E24'; Value = '  +2.26%  ' }
    @{ Cell = 'D25'; Value = '''2.194' }
    @{ Cell = 'E25'; Value = '  +1.02%  ' }
    @{ Cell = 'D26'; Value = '2.141.55' }
    @{ Cell = 'E26'; Value = '  +2.56%  ' }
    @{ Cell = 'D27'; Value = '''21.94' }
    @{ Cell = 'E27'; Value = '  +4.67%  ' }
    @{ Cell = 'D28'; Value = '''162.75' }
    @{ Cell = 'E28'; Value = '  +0.93%  ' }
    @{ Cell = 'D29'; Value = '''2.264' }
    @{ Cell = 'E29'; Value = '  +1.68%  ' }
    @{ Cell = 'D30'; Value = '''128.88' }
    @{ Cell = 'E30'; Value = '  +1.15%  ' }
    @{ Cell = 'D31'; Value = '''1.136' }
    @{ Cell = 'E31'; Value = '  +6.74%  ' }
    @{ Cell = 'D33'; Value = '''5.957' }
    @{ Cell = 'E33'; Value = '  +0.46%  ' }
    @{ Cell = 'D34'; Value = '''3.798' }
    @{ Cell = 'E34'; Value = '  +1.91%  ' }
    @{ Cell = 'D35'; Value = '''0.02451' }
    @{ Cell = 'E35'; Value = '  +0.96%  ' }
    @{ Cell = 'D36'; Value = '''5.315' }
    @{ Cell = 'E36'; Value = '  +2.11%  ' }
    @{ Cell = 'D37'; Value = '''0.06444' }
    @{ Cell = 'E37'; Value = '  +1.56%  ' }
    @{ Cell = 'D38'; Value = '''1.218' }
    @{ Cell = 'E38'; Value = '  +3.87%  ' }
    @{ Cell = 'E39'; Value = '  +0.42%  ' }
    @{ Cell = 'D40'; Value = '''0.6502' }
    @{ Cell = 'E40'; Value = '  +3.18%  ' }
    @{ Cell = 'D41'; Value = '''8.591' }
    @{ Cell = 'E41'; Value = '  +1.10%  ' }
    @{ Cell = 'D42'; Value = '''11.46' }
    @{ Cell = 'E42'; Value = '  +1.58%  ' }
    @{ Cell = 'D43'; Value = '''1.213' }
    @{ Cell = 'E43'; Value = '  +0.81%  ' }
    @{ Cell = 'D44'; Value = '''13.42' }
    @{ Cell = 'E44'; Value = '  +3.36%  ' }
    @{ Cell = 'D45'; Value = '''0.6053' }
    @{ Cell = 'E45'; Value = '  +2.61%  ' }
    @{ Cell = 'D46'; Value = '''2.181' }
    @{ Cell = 'E46'; Value = '  +8.87%  ' }
    @{ Cell = 'D47'; Value = '''3.625' }
    @{ Cell = 'E47'; Value = '  -0.28%  ' }
    @{ Cell = 'D48'; Value = '''122.73' }
    @{ Cell = 'E48'; Value = '  +0.12%  ' }
    @{ Cell = 'E49'; Value = '  +0.80%  ' }
    @{ Cell = 'D50'; Value = '''1.133' }
    @{ Cell = 'E50'; Value = '  +1.24%  ' }
    @{ Cell = 'D51'; Value = '''78.21' }
    @{ Cell = 'E51'; Value = '  +1.67%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
